$d = $word.ActiveDocument
$wNs = 'http://schemas.openxmlformats.org/wordprocessingml/2006/main'

function Find-ParagraphIndex($doc, $needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs.Item($i).Range.Text -like "*$needle*") {
            return $i
        }
    }
    return -1
}

# --- Hunk 1: Heading "Wette:" -> "Frage:" --------------------------------
# Only the run holding the word "Wette" changes text; the trailing ":" run
# (and its own rPr/rsid) is left completely untouched.
$d.Content.Find.Execute("Wette", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Frage", 2) | Out-Null

# --- Hunk 2: "Add.Information: 250" -> spell-checked run split ------------
# The tab stays in its original run; "Add.Information" becomes its own run
# wrapped in proofErr spellStart/spellEnd, and ": 250" becomes a separate run.
$idx = Find-ParagraphIndex $d "Add.Information"
$p = $d.Paragraphs.Item($idx)
$full = $p.Range
$target = $d.Range($full.Start, $full.End - 1)
$xml = '<w:p xmlns:w="' + $wNs + '">' + `
       '<w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:tab/></w:r>' + `
       '<w:proofErr w:type="spellStart"/>' + `
       '<w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t>Add.Information</w:t></w:r>' + `
       '<w:proofErr w:type="spellEnd"/>' + `
       '<w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t>: 250</w:t></w:r>' + `
       '</w:p>'
$target.InsertXML($xml) | Out-Null

# --- Hunk 3: Heading "Wettgruppe:" -> "Grupppe:" (spell-checked split) ----
$idx = Find-ParagraphIndex $d "Wettgruppe"
$p = $d.Paragraphs.Item($idx)
$full = $p.Range
$target = $d.Range($full.Start, $full.End - 1)
$xml = '<w:p xmlns:w="' + $wNs + '">' + `
       '<w:proofErr w:type="spellStart"/>' + `
       '<w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t>Grupppe</w:t></w:r>' + `
       '<w:proofErr w:type="spellEnd"/>' + `
       '<w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t>:</w:t></w:r>' + `
       '</w:p>'
$target.InsertXML($xml) | Out-Null

# --- Hunk 4: "passwordHash: 50" / "nutzername: 25" -------------------------
# passwordHash's value becomes 100 and gets a spell-check split; a brand new
# "phoneNumberHash: 100" paragraph is inserted right after it; "nutzername"
# gets a spell+grammar-checked split and loses its _GoBack bookmark.
$idxPwd = Find-ParagraphIndex $d "passwordHash"
$idxUser = Find-ParagraphIndex $d "nutzername"
$pPwd = $d.Paragraphs.Item($idxPwd)
$pUser = $d.Paragraphs.Item($idxUser)
$target = $d.Range($pPwd.Range.Start, $pUser.Range.End - 1)
$xml = '<w:p xmlns:w="' + $wNs + '">' + `
       '<w:pPr><w:rPr><w:lang w:val="de-DE"/></w:rPr></w:pPr>' + `
       '<w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:tab/></w:r>' + `
       '<w:proofErr w:type="spellStart"/>' + `
       '<w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t>passwordHash</w:t></w:r>' + `
       '<w:proofErr w:type="spellEnd"/>' + `
       '<w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t>: 100</w:t></w:r>' + `
       '</w:p>' + `
       '<w:p xmlns:w="' + $wNs + '">' + `
       '<w:pPr><w:rPr><w:lang w:val="de-DE"/></w:rPr></w:pPr>' + `
       '<w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:tab/></w:r>' + `
       '<w:proofErr w:type="spellStart"/>' + `
       '<w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t>phoneNumberHash</w:t></w:r>' + `
       '<w:proofErr w:type="spellEnd"/>' + `
       '<w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t>: 100</w:t></w:r>' + `
       '</w:p>' + `
       '<w:p xmlns:w="' + $wNs + '">' + `
       '<w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:tab/></w:r>' + `
       '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/>' + `
       '<w:r><w:t>nutzername</w:t></w:r>' + `
       '<w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/>' + `
       '<w:r><w:t>: 25</w:t></w:r>' + `
       '</w:p>'
$target.InsertXML($xml) | Out-Null
